# Week8_Capstone.docx edit script
# 1. "Common Packages" heading -> "Common Themes"
# 2. Expand the Jetpack paragraph with two more sentences
# 3. Add several new paragraphs (ProGuard, Kotlin, Metaprogramming,
#    "Exploring the Manifest" heading, AndroidManifest, Services/receivers,
#    attack-surface / graph paragraph) plus a trailing near-empty paragraph.
# 4. Re-seat the _GoBack bookmark at the new end-of-edit location.

$d = $word.ActiveDocument

# --- 1. "Packages" -> "Themes" in the "Common Packages" Heading2 ---------
$findRng = $d.Content
$null = $findRng.Find.Execute("Packages", $false, $false, $false, $false, $false, $true, 1, $false, $null, 0)
$findRng.Text = "Themes"

# --- 2. Append two sentences to the end of the Jetpack paragraph ---------
# Locate the paragraph that contains the Jetpack sentence so we don't rely
# on hard-coded character offsets.
$jetRng = $d.Content
$null = $jetRng.Find.Execute("accelerates the time to develop new applications. ", $false, $false, $false, $false, $false, $true, 1, $false, $null, 0)
$jetRng.Expand(4)
$insertPoint = $d.Range($jetRng.End - 1, $jetRng.End - 1)
$insertPoint.InsertBefore("Developers can use it from scenarios ranging from interacting with SQLite to animating transitions.")

# --- 3. Insert the new paragraphs right after the Jetpack paragraph ------
# Re-fetch the paragraph range since it grew in step 2.
$jetRng2 = $d.Content
$null = $jetRng2.Find.Execute("accelerates the time to develop new applications. ", $false, $false, $false, $false, $false, $true, 1, $false, $null, 0)
$jetRng2.Expand(4)
$afterJet = $d.Range($jetRng2.End, $jetRng2.End)

$newXml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:tab/></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ProGuard</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> is an open source tool for obfuscation and code reduction that is used by 23.2%</w:t></w:r><w:r><w:t xml:space="preserve"> of sampled applications. These projects are easily identified as all private methods are reduced to </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">single character names. A probable reason for obfuscating open source code is to reduce the final binary size </w:t></w:r><w:r><w:t xml:space="preserve">and </w:t></w:r><w:r><w:t>improve initial download latency.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:tab/><w:t xml:space="preserve">Kotlin is a new language that can be compiled into JavaScript, Java, or into the LLVM compiler framework. </w:t></w:r><w:r><w:t>Only</w:t></w:r><w:r><w:t xml:space="preserve"> 8</w:t></w:r><w:r><w:t>.7</w:t></w:r><w:r><w:t>% of the sampled applications were written in this modern language.</w:t></w:r><w:r><w:t xml:space="preserve"> This was surprisingly low given the reduced barrier to entry and marketing push from Google.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:tab/><w:t xml:space="preserve">Metaprogramming was another theme </w:t></w:r><w:r><w:t xml:space="preserve">that was exposed in many of the common packages such as retrofit, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>butterknife</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, and dagger</w:t></w:r><w:r><w:t>, being used by 10% of applications</w:t></w:r><w:r><w:t>. Traditionally Java has uses reflection however many of the mobile applications leverage Gradle plugins to move the runtime type analysis to compile time.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>Exploring the Manifest</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:tab/><w:t>AndroidManifest.xml describes the components of an application and how they are permitted to interact with the system. The primary entities are activities, providers, services, and receivers. An activity defines the UI behaviors; providers share application content; services are long running background code; and receivers are woken up to handle events (called Intents).</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:tab/><w:t>Services and receivers create an interesting attack vector as they have privileges and do not require user input. For example, a malicious application could send an exploit inside of an Intent object and execute code in the context of another application. Android&#8217;s platform mitigates this scenario by exposing permission</w:t></w:r><w:r><w:t>s on exported application entities</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:tab/></w:r><w:r><w:t xml:space="preserve">To understand the scope of the attack surface 1066 APKs were sampled and 799 of them enabled at least one dangerous permission. </w:t></w:r><w:r><w:t xml:space="preserve">The manifests which use dangerous permissions were parsed into a graph and Gremlin queries executed against it. The graph was able to further reduce </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">the search space to 354 potential targets. </w:t></w:r><w:r><w:t xml:space="preserve">The package </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>at.bitfire.davdroid.apk</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> was selected from the list and confirmed to expose a public service without security policy.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:tab/></w:r></w:p>
"@

$afterJet.InsertXML($newXml)

# --- 4. Re-seat the _GoBack bookmark --------------------------------------
# It should now sit right after "...without security policy." and before
# the final (tab-only) paragraph.
$bmRng = $d.Content
$null = $bmRng.Find.Execute("without security policy.", $false, $false, $false, $false, $false, $true, 1, $false, $null, 0)
$bmPoint = $d.Range($bmRng.End, $bmRng.End)
$d.Bookmarks.Add("_GoBack", $bmPoint)
